$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing data runs from row 2 (date 44075) to row 328 (date 44402).
# Append new rows 329-343 with dates 44403-44417 (26 Jul 2021 - 9 Aug 2021),
# matching the formatting/style of the preceding rows (style index "2" on col A).

$startRow = 329
$startDate = 44403
$endDate = 44417

$r = $startRow
for ($d = $startDate; $d -le $endDate; $d++) {
    $ws.Cells.Item($r, 1).Value = $d
    $ws.Cells.Item($r - 1, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)  # xlPasteFormats

    $ws.Cells.Item($r, 2).Value = 0
    $ws.Cells.Item($r, 3).Value = 0
    $ws.Cells.Item($r, 4).Value = 0

    $r = $r + 1
}
